$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 157.7984646666667
$ws.Range("H2").Value = 473.395394
$ws.Range("I2").Value = 0.341075365555871
$ws.Range("J2").Value = 0.3410753655558709
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.408252
$ws.Range("N2").Value = 19.224756
$ws.Range("O2").Value = 0.8583439096634812
$ws.Range("P2").Value = 0.8583439096634812
$ws.Range("Q2").Value = 1011.212326797096
$ws.Range("R2").Value = 9100.910941173865
$ws.Range("S2").Value = 0.2927599627611273
$ws.Range("T2").Value = 0.2927599627611273
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 157.7984646666667
$ws.Range("H3").Value = 473.395394
$ws.Range("I3").Value = 0.341075365555871
$ws.Range("J3").Value = 0.3410753655558709
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6824433333333334
$ws.Range("N3").Value = 2.04733
$ws.Range("O3").Value = 0.09140887075868921
$ws.Range("P3").Value = 0.09140887075868921
$ws.Range("Q3").Value = 107.6885102220022
$ws.Range("R3").Value = 969.19659199802
$ws.Range("S3").Value = 0.03117731400906929
$ws.Range("T3").Value = 0.03117731400906928
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 157.7984646666667
$ws.Range("H4").Value = 473.395394
$ws.Range("I4").Value = 0.341075365555871
$ws.Range("J4").Value = 0.3410753655558709
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3751373333333333
$ws.Range("N4").Value = 1.125412
$ws.Range("O4").Value = 0.05024721957782962
$ws.Range("P4").Value = 0.05024721957782963
$ws.Range("Q4").Value = 59.19609523914755
$ws.Range("R4").Value = 532.764857152328
$ws.Range("S4").Value = 0.01713808878567435
$ws.Range("T4").Value = 0.01713808878567435
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 127.5109433333333
$ws.Range("H5").Value = 382.53283
$ws.Range("I5").Value = 0.2756100428585324
$ws.Range("J5").Value = 0.2756100428585324
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.408252
$ws.Range("N5").Value = 19.224756
$ws.Range("O5").Value = 0.8583439096634812
$ws.Range("P5").Value = 0.8583439096634812
$ws.Range("Q5").Value = 817.1222576377199
$ws.Range("R5").Value = 7354.10031873948
$ws.Range("S5").Value = 0.2365682017297123
$ws.Range("T5").Value = 0.2365682017297123
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 127.5109433333333
$ws.Range("H6").Value = 382.53283
$ws.Range("I6").Value = 0.2756100428585324
$ws.Range("J6").Value = 0.2756100428585324
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6824433333333334
$ws.Range("N6").Value = 2.04733
$ws.Range("O6").Value = 0.09140887075868921
$ws.Range("P6").Value = 0.09140887075868921
$ws.Range("Q6").Value = 87.01899320487779
$ws.Range("R6").Value = 783.1709388439
$ws.Range("S6").Value = 0.02519320278745238
$ws.Range("T6").Value = 0.02519320278745238
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 127.5109433333333
$ws.Range("H7").Value = 382.53283
$ws.Range("I7").Value = 0.2756100428585324
$ws.Range("J7").Value = 0.2756100428585324
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3751373333333333
$ws.Range("N7").Value = 1.125412
$ws.Range("O7").Value = 0.05024721957782962
$ws.Range("P7").Value = 0.05024721957782963
$ws.Range("Q7").Value = 47.83411525288444
$ws.Range("R7").Value = 430.5070372759599
$ws.Range("S7").Value = 0.01384863834136771
$ws.Range("T7").Value = 0.01384863834136771
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 86.127454
$ws.Range("H8").Value = 258.382362
$ws.Range("I8").Value = 0.1861612083457225
$ws.Range("J8").Value = 0.1861612083457225
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.408252
$ws.Range("N8").Value = 19.224756
$ws.Range("O8").Value = 0.8583439096634812
$ws.Range("P8").Value = 0.8583439096634812
$ws.Range("Q8").Value = 551.9264293504081
$ws.Range("R8").Value = 4967.337864153672
$ws.Range("S8").Value = 0.1597903393991453
$ws.Range("T8").Value = 0.1597903393991453
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 86.127454
$ws.Range("H9").Value = 258.382362
$ws.Range("I9").Value = 0.1861612083457225
$ws.Range("J9").Value = 0.1861612083457225
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6824433333333334
$ws.Range("N9").Value = 2.04733
$ws.Range("O9").Value = 0.09140887075868921
$ws.Range("P9").Value = 0.09140887075868921
$ws.Range("Q9").Value = 58.77710679927334
$ws.Range("R9").Value = 528.99396119346
$ws.Range("S9").Value = 0.01701678583395556
$ws.Range("T9").Value = 0.01701678583395556
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 86.127454
$ws.Range("H10").Value = 258.382362
$ws.Range("I10").Value = 0.1861612083457225
$ws.Range("J10").Value = 0.1861612083457225
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3751373333333333
$ws.Range("N10").Value = 1.125412
$ws.Range("O10").Value = 0.05024721957782962
$ws.Range("P10").Value = 0.05024721957782963
$ws.Range("Q10").Value = 32.30962342034933
$ws.Range("R10").Value = 290.7866107831439
$ws.Range("S10").Value = 0.009354083112621606
$ws.Range("T10").Value = 0.009354083112621606
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 91.212982
$ws.Range("H11").Value = 273.638946
$ws.Range("I11").Value = 0.1971533832398742
$ws.Range("J11").Value = 0.1971533832398741
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.408252
$ws.Range("N11").Value = 19.224756
$ws.Range("O11").Value = 0.8583439096634812
$ws.Range("P11").Value = 0.8583439096634812
$ws.Range("Q11").Value = 584.515774327464
$ws.Range("R11").Value = 5260.641968947175
$ws.Range("S11").Value = 0.1692254057734962
$ws.Range("T11").Value = 0.1692254057734962
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 91.212982
$ws.Range("H12").Value = 273.638946
$ws.Range("I12").Value = 0.1971533832398742
$ws.Range("J12").Value = 0.1971533832398741
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6824433333333334
$ws.Range("N12").Value = 2.04733
$ws.Range("O12").Value = 0.09140887075868921
$ws.Range("P12").Value = 0.09140887075868921
$ws.Range("Q12").Value = 62.24769147935334
$ws.Range("R12").Value = 560.22922331418
$ws.Range("S12").Value = 0.01802156812821198
$ws.Range("T12").Value = 0.01802156812821198
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 91.212982
$ws.Range("H13").Value = 273.638946
$ws.Range("I13").Value = 0.1971533832398742
$ws.Range("J13").Value = 0.1971533832398741
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.3751373333333333
$ws.Range("N13").Value = 1.125412
$ws.Range("O13").Value = 0.05024721957782962
$ws.Range("P13").Value = 0.05024721957782963
$ws.Range("Q13").Value = 34.21739483286132
$ws.Range("R13").Value = 307.956553495752
$ws.Range("S13").Value = 0.009906409338165952
$ws.Range("T13").Value = 0.009906409338165952
